$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G
$data = @(
    @(-0.3309730362341788, 0.6874379636966282, 1.909552438266325, 1.381865564469397, 1.357907196056604, 42),
    @(0.03815854055786114, 1.119444934519008, 3.277134428859379, 1.810285731275419, 1.832367406801354, 41),
    @(-0.2513998558810283, 0.8541694100438215, 1.593836722990716, 1.262472464250495, 1.252949273720461, 40),
    @(-0.07458524723922935, 0.8386825321727622, 1.828534390886613, 1.352233112627631, 1.367824683282241, 39),
    @(-0.2016351845399217, 0.8634258771531403, 1.6349713752315, 1.278659992035217, 1.279610961236435, 38),
    @(-0.05768949108510493, 1.131474238540056, 2.649968512540811, 1.627872388285031, 1.649290218556056, 37),
    @(-0.2756624231401172, 0.8747661734718721, 1.641880542173012, 1.281358865491246, 1.269106184267327, 36),
    @(-0.04407887464051687, 0.8715491677183831, 1.335121592579382, 1.155474617886253, 1.171490410553548, 35),
    @(-0.2588473913255688, 0.7912842664960333, 1.430082320060018, 1.195860493561025, 1.185067851710109, 34),
    @(0.02168306473352078, 1.183415881580869, 2.744833036578304, 1.656753764618721, 1.682297304676594, 33)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j  # B=2 ... G=7
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
